# formatação cep e telefone
# Converts the numeric CEP (E) and Telefone (F) columns into formatted
# text strings, e.g. 54768122 -> "54768-122" and 12345678901 -> "(12) 3456-678901"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Arthur
$ws.Range("E2").Value = "54768-122"
$ws.Range("F2").Value = "(12) 3456-678901"

# Row 3 - Lucas
$ws.Range("E3").Value = "54768-122"
$ws.Range("F3").Value = "(99) 3456-678901"

# Row 4 - Lucas
$ws.Range("E4").Value = "54768-122"
$ws.Range("F4").Value = "(99) 3456-678901"

# Row 5 - Bola
$ws.Range("E5").Value = "54768-122"
$ws.Range("F5").Value = "(99) 3456-678989"

# Row 6 - Cabeça
$ws.Range("E6").Value = "54768-122"
$ws.Range("F6").Value = "(99) 3456-678989"
